# Helper: write a value that *looks* numeric (leading zeros, decimals as
# text, etc.) while keeping it a genuine text cell with default ("Normal")
# styling - mirrors how the source workbook stores e.g. "014339" or "3.28"
# as text rather than as a number.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right before "总计", using "2021-Q4"
#    as a formatting template (same column layout: A=index, B..H =
#    fund code / name / scale / total position / position % / value /
#    rank).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Re-fetch sheet references by name - adding a sheet can shift/stale any
# previously-held worksheet references in this object model.
$template = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Item("2022-Q1")

# Copy the header + 6 data row formatting template (A1:H7) from 2021-Q4.
$template.Range("A1:H7").Copy()
$ws.Range("A1:H7").PasteSpecial(-4122)

# Header row
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# Data rows: idx, code, name, scale, totalPos, posPct, value, rank
$rows = @(
    @(0, "014339", "长江智能制造混合A",     "3.28", "21.63", "1.46", "0.0479", 6),
    @(1, "009128", "明亚价值长青混合A",     "0.38", "49.48", "4.15", "0.0158", 3),
    @(2, "000892", "九泰天宝灵活配置混合A", "0.07", "90.81", "4.52", "0.0032", 9),
    @(3, "014340", "长江智能制造混合C",     "0.15", "21.63", "1.46", "0.0022", 6),
    @(4, "002028", "九泰天宝灵活配置混合C", "0.00", "90.81", "4.52", "0",      9),
    @(5, "009129", "明亚价值长青混合C",     "0.00", "49.48", "4.15", "0",      3)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Cells.Item($r,1).Value = $row[0]

    Set-TextValue $ws.Cells.Item($r,2) $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    Set-TextValue $ws.Cells.Item($r,4) $row[3]
    Set-TextValue $ws.Cells.Item($r,5) $row[4]
    Set-TextValue $ws.Cells.Item($r,6) $row[5]

    if ($row[6] -eq "0") {
        $ws.Cells.Item($r,7).Value = 0
    } else {
        Set-TextValue $ws.Cells.Item($r,7) $row[6]
    }

    $ws.Cells.Item($r,8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: add a new top data row for 2022-Q1 and
#    shift the existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the sheet by one row, copying the per-column formatting from the
# current last row (row 6) down into the new row 7. Restricted to the
# used A:D range so we don't touch/format the entire 16384-column row.
$total.Range("A6:D6").Copy()
$total.Range("A7:D7").PasteSpecial(-4122)

# Snapshot existing data rows (2..6) before overwriting them.
$oldData = @()
for ($r = 2; $r -le 6; $r++) {
    $oldData += ,@($total.Cells.Item($r,2).Value2, $total.Cells.Item($r,3).Value2, $total.Cells.Item($r,4).Value2)
}

# Shift rows 2..6 down to 3..7 (old index i -> new row i+3), refreshing
# the running index in column A.
for ($i = 4; $i -ge 0; $i--) {
    $row = $i + 3
    $total.Cells.Item($row,1).Value = $i + 1
    $total.Cells.Item($row,2).Value = $oldData[$i][0]
    $total.Cells.Item($row,3).Value = $oldData[$i][1]
    $total.Cells.Item($row,4).Value = $oldData[$i][2]
}

# Write the new 2022-Q1 summary row at row 2.
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 6
$total.Cells.Item(2,4).Value = 0.07
